$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-11
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
